$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-11-26 18:22:19"

for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
